# feat(commands): adapt xlspython to the new MVC architecture
#
# A new "raw code" column is inserted right before the existing
# "Colonne de(s) maximum(s)" column (F). The old column F (and all of its
# data/header) slides one column to the right and becomes column G. The
# freshly inserted column F is then populated, for the rows that already
# carry a code in column E, with that same code value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new blank column at F - this shifts the previous F
# column (header + the "a"/"af"/"g"/"ff" values) one place over to G,
# carrying its values and formatting along with it.
$ws.Columns("F").Insert()

# Column E already holds the condition code for a handful of rows
# (prime_congruent / probe_incongruent / probe_congruent / prime_neutre).
# Mirror those same codes into the newly created column F, matching the
# look of the surrounding cells (same style as the neighbouring column).
$codeRows = @(2, 3, 6, 9)
foreach ($r in $codeRows) {
    $srcCell = $ws.Range("E" + $r)
    $dstCell = $ws.Range("F" + $r)
    $dstCell.Value = $srcCell.Value2
    $dstCell.Style = $ws.Range("G" + $r).Style
}
